$d = $word.ActiveDocument
$old = "setStyle(org.obeonetwork.m2doc.element.MParagraph,java.lang.String) with arguments [org.obeonetwork.m2doc.element.impl.MParagraphImpl@55fdf7f9, NotExisting] failed:`n`tno text style NotExisting`njava.lang.IllegalArgumentException: no text style NotExisting`n`tat org.obeonetwork.m2doc.services.PaginationServices.setStyle(PaginationServices.java:366)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)`n`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:208)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)`n`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:109)`n`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:591)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1556)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:301)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:286)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)`n`tat sun.reflect.GeneratedMethodAccessor4.invoke(Unknown Source)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:59)`n`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)`n`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:56)`n`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.BlockJUnit4ClassRunner`$1.evaluate(BlockJUnit4ClassRunner.java:100)`n`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:366)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:103)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:63)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$4.run(ParentRunner.java:331)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:79)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:329)`n`tat org.junit.runners.ParentRunner.access`$100(ParentRunner.java:66)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:293)`n`tat org.junit.runners.ParentRunner`$3.evaluate(ParentRunner.java:306)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:413)`n`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)`n"
$new = "setStyle(org.obeonetwork.m2doc.element.MParagraph,java.lang.String) with arguments [org.obeonetwork.m2doc.element.impl.MParagraphImpl@106b014e, NotExisting] failed:`n`tno text style NotExisting`njava.lang.IllegalArgumentException: no text style NotExisting`n`tat org.obeonetwork.m2doc.services.PaginationServices.setStyle(PaginationServices.java:366)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:162)`n`tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:135)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:172)`n`tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callOrApply(EvaluationServices.java:208)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:192)`n`tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:119)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)`n`tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:604)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1675)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:314)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)`n`tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)`n`tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1450)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:299)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:536)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:422)`n`tat sun.reflect.GeneratedMethodAccessor6.invoke(Unknown Source)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.junit.runners.model.FrameworkMethod`$1.runReflectiveCall(FrameworkMethod.java:50)`n`tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)`n`tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)`n`tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)`n`tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)`n`tat org.junit.runners.Suite.runChild(Suite.java:27)`n`tat org.junit.runners.ParentRunner`$3.run(ParentRunner.java:290)`n`tat org.junit.runners.ParentRunner`$1.schedule(ParentRunner.java:71)`n`tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)`n`tat org.junit.runners.ParentRunner.access`$000(ParentRunner.java:58)`n`tat org.junit.runners.ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)`n"
$found = $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Host "Found: $found"
